$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab05")

# Fix mojibake (corrupted Latin-1/UTF-8 round trip) characters in the
# "Regional Economic Communities" footnote text (cell A103):
#   Pa>ses -> Países, L>ngua -> Língua, Com>n -> Común
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Small precision correction on L72
$ws.Range("L72").Value = 78.326488430312494

# Data corrections for row 97 ("Africa, Fragile States")
$ws.Range("C97").Value = 842549.25199999998
$ws.Range("D97").Value = 348533.71888524003
$ws.Range("E97").Value = 494015.53311476001
$ws.Range("F97").Value = 48796.272602885401
$ws.Range("G97").Value = 85
$ws.Range("H97").Value = 117932.369893483
$ws.Range("I97").Value = 39
$ws.Range("J97").Value = 5.4141204157620599
$ws.Range("K97").Value = 72.908692094416395
$ws.Range("L97").Value = 78.322812510178395

# Data corrections for row 98 ("ROW, Fragile States")
$ws.Range("C98").Value = 692226.44200000004
$ws.Range("D98").Value = 318611.69368165999
$ws.Range("E98").Value = 373614.74831833999
$ws.Range("F98").Value = 36149.9535672317
$ws.Range("G98").Value = 62
$ws.Range("H98").Value = 129607.049632435
$ws.Range("I98").Value = 33
$ws.Range("J98").Value = 8.7123928268905999
$ws.Range("K98").Value = 47.823885843296999
$ws.Range("L98").Value = 56.5362786701875
